$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 80.42856999999999
$ws.Range("I2").Value = 68.833336
$ws.Range("K2").Value = 68.833336
$ws.Range("M2").Value = 44.166664
$ws.Range("H33").Value = 577.86365
$ws.Range("I33").Value = 276.9375
$ws.Range("J33").Value = 1380.3334
$ws.Range("K33").Value = 276.9375
$ws.Range("L33").Value = 1380.3334
$ws.Range("M33").Value = -47.9375
$ws.Range("N33").Value = -1838.3334
$ws.Range("H70").Value = 11906136
$ws.Range("I70").Value = 847.2
$ws.Range("J70").Value = 22729126
$ws.Range("K70").Value = 2541.6
$ws.Range("L70").Value = 68187378
$ws.Range("M70").Value = -2271.6
$ws.Range("N70").Value = -68187918
$ws.Range("H73").Value = 11906136
$ws.Range("I73").Value = 847.2
$ws.Range("J73").Value = 22729126
$ws.Range("K73").Value = 2541.6
$ws.Range("L73").Value = 68187378
$ws.Range("M73").Value = -1605.6
$ws.Range("N73").Value = -68189250
$ws.Range("H98").Value = 1769.3077
$ws.Range("I98").Value = 1347.6364
$ws.Range("K98").Value = 1347.6364
$ws.Range("M98").Value = 150.3635999999999
$ws.Range("H112").Value = 4219.6665
$ws.Range("I112").Value = 1868
$ws.Range("J112").Value = 5395.5
$ws.Range("K112").Value = 5604
$ws.Range("L112").Value = 16186.5
$ws.Range("M112").Value = -4496
$ws.Range("N112").Value = -18402.5
$ws.Range("H121").Value = 4004.9412
$ws.Range("J121").Value = 4004.9412
$ws.Range("L121").Value = 12014.8236
$ws.Range("N121").Value = -15508.8236
$ws.Range("H122").Value = 1769.3077
$ws.Range("I122").Value = 1347.6364
$ws.Range("K122").Value = 4042.9092
$ws.Range("M122").Value = -1592.9092
$ws.Range("H138").Value = 5075.7
$ws.Range("I138").Value = 1731.1666
$ws.Range("K138").Value = 5193.4998
$ws.Range("M138").Value = -53.4997999999996
$ws.Range("H140").Value = 580747
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 2926.2
$ws.Range("I141").Value = 2926.2
$ws.Range("K141").Value = 8778.599999999999
$ws.Range("M141").Value = -3598.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4254284.5
$ws.Range("I32").Value = 1794520.5
$ws.Range("K32").Value = 1794520.5
$ws.Range("M32").Value = -1794233.5
$ws.Range("H61").Value = 3905.4666
$ws.Range("I61").Value = 3621.6924
$ws.Range("J61").Value = 5750
$ws.Range("K61").Value = 3621.6924
$ws.Range("L61").Value = 5750
$ws.Range("M61").Value = -3409.6924
$ws.Range("N61").Value = -6174
$ws.Range("H74").Value = 3784.3845
$ws.Range("I74").Value = 3717.1904
$ws.Range("K74").Value = 3717.1904
$ws.Range("M74").Value = -2843.1904
$ws.Range("H77").Value = 3784.3845
$ws.Range("I77").Value = 3717.1904
$ws.Range("K77").Value = 18585.952
$ws.Range("M77").Value = -14217.952
$ws.Range("H132").Value = 2487.9375
$ws.Range("I132").Value = 2094.6584
$ws.Range("J132").Value = 4791.4287
$ws.Range("K132").Value = 6283.975199999999
$ws.Range("L132").Value = 14374.2861
$ws.Range("M132").Value = -3753.975199999999
$ws.Range("N132").Value = -19434.2861
$ws.Range("H136").Value = 3905.4666
$ws.Range("I136").Value = 3621.6924
$ws.Range("J136").Value = 5750
$ws.Range("K136").Value = 10865.0772
$ws.Range("L136").Value = 17250
$ws.Range("M136").Value = -8315.0772
$ws.Range("N136").Value = -22350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2232098.2
$ws.Range("J7").Value = 16666.334
$ws.Range("L7").Value = 16666.334
$ws.Range("N7").Value = -16892.334
$ws.Range("H80").Value = 638.6667
$ws.Range("J80").Value = 633.0625
$ws.Range("L80").Value = 633.0625
$ws.Range("N80").Value = -2629.0625
$ws.Range("H83").Value = 638.6667
$ws.Range("J83").Value = 633.0625
$ws.Range("L83").Value = 3165.3125
$ws.Range("N83").Value = -13149.3125
$ws.Range("H107").Value = 1197.4482
$ws.Range("J107").Value = 1347.4
$ws.Range("L107").Value = 1347.4
$ws.Range("N107").Value = -5187.4
$ws.Range("H134").Value = 6495851.5
$ws.Range("I134").Value = 6495851.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 19487554.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -19485019.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3649.5
$ws.Range("I16").Value = 3999
$ws.Range("J16").Value = 3439.8
$ws.Range("K16").Value = 3999
$ws.Range("L16").Value = 3439.8
$ws.Range("M16").Value = -3712
$ws.Range("N16").Value = -4013.8
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81748
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -248736
$ws.Range("H99").Value = 4477.8335
$ws.Range("I99").Value = 2632.6667
$ws.Range("K99").Value = 2632.6667
$ws.Range("M99").Value = -1134.6667
$ws.Range("H105").Value = 1842.3
$ws.Range("I105").Value = 1630.3334
$ws.Range("J105").Value = 3750
$ws.Range("K105").Value = 1630.3334
$ws.Range("L105").Value = 3750
$ws.Range("M105").Value = 116.6666
$ws.Range("N105").Value = -7244
$ws.Range("H113").Value = 3649.5
$ws.Range("I113").Value = 3999
$ws.Range("J113").Value = 3439.8
$ws.Range("K113").Value = 3999
$ws.Range("L113").Value = 3439.8
$ws.Range("M113").Value = -1829
$ws.Range("N113").Value = -7779.8
$ws.Range("H126").Value = 4477.8335
$ws.Range("I126").Value = 2632.6667
$ws.Range("K126").Value = 7898.000100000001
$ws.Range("M126").Value = -5428.000100000001
$ws.Range("H132").Value = 2410.8572
$ws.Range("I132").Value = 2410.8572
$ws.Range("K132").Value = 7232.571599999999
$ws.Range("M132").Value = -4702.571599999999
$ws.Range("H134").Value = 3487.4119
$ws.Range("I134").Value = 3640.182
$ws.Range("K134").Value = 10920.546
$ws.Range("M134").Value = -8385.545999999998
$ws.Range("H141").Value = 360013.22
$ws.Range("J141").Value = 360013.22
$ws.Range("L141").Value = 360013.22
$ws.Range("N141").Value = -370373.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4273.2856
$ws.Range("J39").Value = 4985
$ws.Range("L39").Value = 14955
$ws.Range("N39").Value = -15543
$ws.Range("H75").Value = 799.6667
$ws.Range("I75").Value = 599.5
$ws.Range("J75").Value = 1200
$ws.Range("K75").Value = 1798.5
$ws.Range("L75").Value = 3600
$ws.Range("M75").Value = -800.5
$ws.Range("N75").Value = -5596
$ws.Range("H78").Value = 799.6667
$ws.Range("I78").Value = 599.5
$ws.Range("J78").Value = 1200
$ws.Range("K78").Value = 5395.5
$ws.Range("L78").Value = 10800
$ws.Range("M78").Value = -403.5
$ws.Range("N78").Value = -20784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3050
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3050
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -29984
$ws.Range("H132").Value = 2207.9092
$ws.Range("I132").Value = 1435.875
$ws.Range("J132").Value = 4266.6665
$ws.Range("K132").Value = 4307.625
$ws.Range("L132").Value = 12799.9995
$ws.Range("M132").Value = -1777.625
$ws.Range("N132").Value = -17859.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5663.8213
$ws.Range("I7").Value = 5482.9165
$ws.Range("K7").Value = 5482.9165
$ws.Range("M7").Value = -5370.9165
$ws.Range("H46").Value = 2798.111
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 3054.7144
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 3054.7144
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -3430.7144
$ws.Range("H61").Value = 3575.1177
$ws.Range("I61").Value = 5522.4443
$ws.Range("K61").Value = 5522.4443
$ws.Range("M61").Value = -5320.4443
$ws.Range("H113").Value = 3575.1177
$ws.Range("I113").Value = 5522.4443
$ws.Range("K113").Value = 5522.4443
$ws.Range("M113").Value = -3352.4443
$ws.Range("H126").Value = 5663.8213
$ws.Range("I126").Value = 5482.9165
$ws.Range("K126").Value = 16448.7495
$ws.Range("M126").Value = -13978.7495
$ws.Range("H132").Value = 3300
$ws.Range("I132").Value = 2600
$ws.Range("K132").Value = 7800
$ws.Range("M132").Value = -5270

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2415.25
$ws.Range("I132").Value = 2263.2058
$ws.Range("K132").Value = 6789.617400000001
$ws.Range("M132").Value = -4259.617400000001
$ws.Range("H136").Value = 1846.2222
$ws.Range("I136").Value = 1282.9395
$ws.Range("K136").Value = 3848.8185
$ws.Range("M136").Value = -1298.8185

